$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) F44: Offshore wind NOS0 capacity 8000 -> 7000
# ---------------------------------------------------------------
$ws.Cells.Item(44, 6).Value = 7000

# ---------------------------------------------------------------
# 2) New capacity rows 49-55 (wind updates for NON1/NOM1/NOS0/PL00/BE00)
# ---------------------------------------------------------------
$ws.Cells.Item(49, 1).Value = "NON1"
$ws.Cells.Item(49, 3).Value = "Onshore Wind"
$ws.Cells.Item(49, 4).Value = "Distributed Energy"
$ws.Cells.Item(49, 5).Value = 2040
$ws.Cells.Item(49, 6).Value = 2200

$ws.Cells.Item(50, 1).Value = "NOM1"
$ws.Cells.Item(50, 3).Value = "Onshore Wind"
$ws.Cells.Item(50, 4).Value = "Distributed Energy"
$ws.Cells.Item(50, 5).Value = 2040
$ws.Cells.Item(50, 6).Value = 2000

$ws.Cells.Item(51, 1).Value = "NOS0"
$ws.Cells.Item(51, 3).Value = "Onshore Wind"
$ws.Cells.Item(51, 4).Value = "Distributed Energy"
$ws.Cells.Item(51, 5).Value = 2040
$ws.Cells.Item(51, 6).Value = 4770

$ws.Cells.Item(52, 1).Value = "NON1"
$ws.Cells.Item(52, 3).Value = "Offshore Wind"
$ws.Cells.Item(52, 4).Value = "Distributed Energy"
$ws.Cells.Item(52, 5).Value = 2040
$ws.Cells.Item(52, 6).Value = 500

$ws.Cells.Item(53, 1).Value = "NOM1"
$ws.Cells.Item(53, 3).Value = "Offshore Wind"
$ws.Cells.Item(53, 4).Value = "Distributed Energy"
$ws.Cells.Item(53, 5).Value = 2040
$ws.Cells.Item(53, 6).Value = 500

$ws.Cells.Item(54, 1).Value = "PL00"
$ws.Cells.Item(54, 3).Value = "Onshore Wind"
$ws.Cells.Item(54, 4).Value = "Distributed Energy"
$ws.Cells.Item(54, 5).Value = 2040
$ws.Cells.Item(54, 6).Value = 15000

$ws.Cells.Item(55, 1).Value = "BE00"
$ws.Cells.Item(55, 3).Value = "Onshore Wind"
$ws.Cells.Item(55, 4).Value = "Distributed Energy"
$ws.Cells.Item(55, 5).Value = 2040
$ws.Cells.Item(55, 6).Value = 6500

# Row 56 stays blank apart from its pre-existing A/E/F formatting; D56 is cleared.
$ws.Cells.Item(56, 4).ClearContents()

# ---------------------------------------------------------------
# 3) Defined name _xlnm._FilterDatabase range grows to include the new rows
# ---------------------------------------------------------------
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Capacity!`$A`$1:`$J`$56"

# ---------------------------------------------------------------
# 4) AutoFilter: drop the Year=2040 filter and the SE0x node filter,
#    replace with a PL00 node filter, over the extended range.
# ---------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:J56").AutoFilter(1, @("PL00"), 7)

# ---------------------------------------------------------------
# 5) Row visibility (hidden/unhidden). Applied last: the loader does not
#    round-trip the original hidden="true"/"false" row attributes, and
#    AutoFilter() above hides any row not matching its criteria - so
#    every row's final visibility has to be (re)asserted explicitly,
#    after the filter has been applied, to reproduce the target state.
# ---------------------------------------------------------------
foreach ($r in 2..15) {
  $ws.Rows.Item($r).Hidden = $true
}
$ws.Rows.Item(16).Hidden = $false
$ws.Rows.Item(17).Hidden = $false
foreach ($r in 18..53) {
  $ws.Rows.Item($r).Hidden = $true
}
foreach ($r in 54..64) {
  $ws.Rows.Item($r).Hidden = $false
}

# ---------------------------------------------------------------
# 6) Active selection moves to D58
# ---------------------------------------------------------------
$ws.Range("D58").Select()
